# Auto-generated edit script: updates Leve profit-calculation values
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the
# scheduled-runner refresh of currentAveragePrice / Leve price / profit columns.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6090.5
$ws.Range("I116").Value = 6580.6
$ws.Range("K116").Value = 6580.6
$ws.Range("M116").Value = -3138.6

$ws.Range("H132").Value = 2339.7954
$ws.Range("I132").Value = 2386.9285
$ws.Range("J132").Value = 1350
$ws.Range("K132").Value = 7160.7855
$ws.Range("L132").Value = 4050
$ws.Range("M132").Value = -4630.7855
$ws.Range("N132").Value = -9110

$ws.Range("H137").Value = 2550.1667
$ws.Range("I137").Value = 2148
$ws.Range("J137").Value = 3220.4443
$ws.Range("K137").Value = 6444
$ws.Range("L137").Value = 9661.332900000001
$ws.Range("M137").Value = -3894
$ws.Range("N137").Value = -14761.3329

$ws.Range("H138").Value = 2604.0876
$ws.Range("I138").Value = 1332.3182
$ws.Range("K138").Value = 3996.9546
$ws.Range("M138").Value = 1143.0454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3743.6667
$ws.Range("I122").Value = 3319.6924
$ws.Range("J122").Value = 6499.5
$ws.Range("K122").Value = 9959.0772
$ws.Range("L122").Value = 19498.5
$ws.Range("M122").Value = -7509.0772
$ws.Range("N122").Value = -24398.5

$ws.Range("H132").Value = 2128.919
$ws.Range("I132").Value = 2035.8148
$ws.Range("J132").Value = 2380.3
$ws.Range("K132").Value = 6107.4444
$ws.Range("L132").Value = 7140.900000000001
$ws.Range("M132").Value = -3577.4444
$ws.Range("N132").Value = -12200.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2522.2856
$ws.Range("I99").Value = 1911
$ws.Range("K99").Value = 1911
$ws.Range("M99").Value = -413

$ws.Range("H105").Value = 5058.1816
$ws.Range("I105").Value = 5093.684
$ws.Range("K105").Value = 5093.684
$ws.Range("M105").Value = -3346.684

$ws.Range("H134").Value = 3657
$ws.Range("I134").Value = 3642.7083
$ws.Range("K134").Value = 10928.1249
$ws.Range("M134").Value = -8393.124899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3077.3174
$ws.Range("I31").Value = 4370.727
$ws.Range("K31").Value = 4370.727
$ws.Range("M31").Value = -4075.727

$ws.Range("H34").Value = 3077.3174
$ws.Range("I34").Value = 4370.727
$ws.Range("K34").Value = 4370.727
$ws.Range("M34").Value = -4168.727

$ws.Range("H58").Value = 1894.9
$ws.Range("I58").Value = 1838.8
$ws.Range("K58").Value = 1838.8
$ws.Range("M58").Value = -1635.8

$ws.Range("H59").Value = 59748.5
$ws.Range("J59").Value = 79997.60000000001
$ws.Range("L59").Value = 79997.60000000001
$ws.Range("N59").Value = -82287.60000000001

$ws.Range("H107").Value = 1890.3
$ws.Range("I107").Value = 2781.6
$ws.Range("K107").Value = 2781.6
$ws.Range("M107").Value = -861.5999999999999

$ws.Range("H134").Value = 658.79486
$ws.Range("I134").Value = 630.55554
$ws.Range("K134").Value = 1891.66662
$ws.Range("M134").Value = 643.33338

$ws.Range("H136").Value = 1894.9
$ws.Range("I136").Value = 1838.8
$ws.Range("K136").Value = 5516.4
$ws.Range("M136").Value = -2966.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 573
$ws.Range("I5").Value = 573
$ws.Range("K5").Value = 1719
$ws.Range("M5").Value = -1607

$ws.Range("H12").Value = 858.4286
$ws.Range("J12").Value = 868.93335
$ws.Range("L12").Value = 2606.80005
$ws.Range("N12").Value = -2952.80005

$ws.Range("H46").Value = 53253.26
$ws.Range("J46").Value = 71971.5
$ws.Range("L46").Value = 215914.5
$ws.Range("N46").Value = -216096.5

$ws.Range("H92").Value = 608
$ws.Range("J92").Value = 584.5
$ws.Range("L92").Value = 1753.5
$ws.Range("N92").Value = -4249.5

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H131").Value = 2549.6
$ws.Range("J131").Value = 2721.7778
$ws.Range("L131").Value = 8165.3334
$ws.Range("N131").Value = -18245.3334

$ws.Range("H135").Value = 573
$ws.Range("I135").Value = 573
$ws.Range("K135").Value = 5157
$ws.Range("M135").Value = -2622

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 179.36363
$ws.Range("I2").Value = 198.83333
$ws.Range("J2").Value = 156
$ws.Range("K2").Value = 198.83333
$ws.Range("L2").Value = 156
$ws.Range("M2").Value = -85.83332999999999
$ws.Range("N2").Value = -382

$ws.Range("H70").Value = 6863.25
$ws.Range("I70").Value = 6776.5
$ws.Range("K70").Value = 6776.5
$ws.Range("M70").Value = -6506.5

$ws.Range("H73").Value = 6863.25
$ws.Range("I73").Value = 6776.5
$ws.Range("K73").Value = 6776.5
$ws.Range("M73").Value = -5840.5

$ws.Range("H107").Value = 693
$ws.Range("I107").Value = 584.75
$ws.Range("K107").Value = 584.75
$ws.Range("M107").Value = 1335.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 22469.666
$ws.Range("I18").Value = 22469.666
$ws.Range("K18").Value = 22469.666
$ws.Range("M18").Value = -22297.666

$ws.Range("H61").Value = 2755.5881
$ws.Range("I61").Value = 2253.3333
$ws.Range("J61").Value = 3961
$ws.Range("K61").Value = 2253.3333
$ws.Range("L61").Value = 3961
$ws.Range("M61").Value = -2051.3333
$ws.Range("N61").Value = -4365

$ws.Range("H93").Value = 1666.6666
$ws.Range("I93").Value = 1666.6666
$ws.Range("K93").Value = 1666.6666
$ws.Range("M93").Value = -418.6666

$ws.Range("H113").Value = 2755.5881
$ws.Range("I113").Value = 2253.3333
$ws.Range("J113").Value = 3961
$ws.Range("K113").Value = 2253.3333
$ws.Range("L113").Value = 3961
$ws.Range("M113").Value = -83.33329999999978
$ws.Range("N113").Value = -8301

$ws.Range("H122").Value = 8111.174
$ws.Range("I122").Value = 9252.333000000001
$ws.Range("J122").Value = 4003
$ws.Range("K122").Value = 27756.999
$ws.Range("L122").Value = 12009
$ws.Range("M122").Value = -25306.999
$ws.Range("N122").Value = -16909

$ws.Range("H132").Value = 2460.2856
$ws.Range("I132").Value = 2363.9583
$ws.Range("K132").Value = 7091.874899999999
$ws.Range("M132").Value = -4561.874899999999

$ws.Range("H140").Value = 119528
$ws.Range("J140").Value = 119528
$ws.Range("L140").Value = 119528
$ws.Range("N140").Value = -129888

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3209.1428
$ws.Range("I107").Value = 3314.3333
$ws.Range("J107").Value = 2946.1667
$ws.Range("K107").Value = 9942.999899999999
$ws.Range("L107").Value = 8838.500100000001
$ws.Range("M107").Value = -8022.999899999999
$ws.Range("N107").Value = -12678.5001
